$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A couple of "Price" cells get a new value that still parses as a plain
# number (e.g. "14.80"). Left alone, Excel auto-converts such text into a
# numeric literal and silently drops the insignificant trailing zero,
# which would corrupt the displayed text. Mark those specific cells as
# Text before writing so the literal string is preserved exactly.
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"

$ws.Range('D2').Value = '26.694.80'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '1.794.18'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '309.13'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '0.4457'
$ws.Range('E7').Value = '  +5.31%  '
$ws.Range('D8').Value = '0.3656'
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('D9').Value = '0.07294'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '0.8549'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '1.977.55'
$ws.Range('E11').Value = '  +8.30%  '
$ws.Range('D12').Value = '20.56'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').Value = '6.598'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.07072'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = '91.94'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '5.265'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '0.000008668'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '14.80'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').Value = '26.762.02'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').Value = '5.138'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').Value = '10.75'
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('D24').Value = '1.984'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '151.90'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '18.42'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.168'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('D28').Value = '5.175'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('D29').Value = '116.28'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').Value = '0.08786'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').Value = '0.7401'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').Value = '1.152'
$ws.Range('E32').Value = '  -3.23%  '
$ws.Range('D33').Value = '2.935'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('D34').Value = '4.423'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').Value = '1.083'
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('D37').Value = '0.01955'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').Value = '0.05168'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').Value = '0.5268'
$ws.Range('E39').Value = '  +4.59%  '
$ws.Range('D40').Value = '2.845'
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('D41').Value = '7.019'
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('D42').Value = '0.1676'
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').Value = '0.5085'
$ws.Range('E43').Value = '  +7.31%  '
$ws.Range('D44').Value = '8.371'
$ws.Range('E44').Value = '  -3.26%  '
$ws.Range('D45').Value = '10.44'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').Value = '1.957'
$ws.Range('E46').Value = '  +4.32%  '
$ws.Range('D47').Value = '105.44'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Value = '1.655'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').Value = '0.06297'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('D51').Value = '0.9122'
$ws.Range('E51').Value = '  -0.01%  '
